$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("G2").Value = 3.3
$ws.Range("I2").Value = 2.4
$ws.Range("K2").Value = 1.95
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("O2").Value = 1.4
$ws.Range("P2").Value = 2.75
$ws.Range("Q2").Value = 1.78
$ws.Range("R2").Value = 2.1
$ws.Range("S2").Value = 2.35
$ws.Range("T2").Value = 1.57
$ws.Range("W2").Value = 4.33
$ws.Range("X2").Value = 1.2
$ws.Range("AC2").Value = 8.5
$ws.Range("AO2").Value = 10

# Row 4 updates
$ws.Range("G4").Value = 1.75
$ws.Range("I4").Value = 4.2
$ws.Range("J4").Value = 2.5
$ws.Range("K4").Value = 2.1
$ws.Range("L4").Value = 5
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 10
$ws.Range("S4").Value = 2.15
$ws.Range("T4").Value = 1.67
$ws.Range("Y4").Value = 1.44
$ws.Range("Z4").Value = 2.63
$ws.Range("AA4").Value = 2
$ws.Range("AB4").Value = 1.73
$ws.Range("AK4").Value = 19
$ws.Range("AO4").Value = 21
$ws.Range("AS4").Value = 41

# Row 5 updates
$ws.Range("N5").Value = 9
